$d = $word.ActiveDocument

# Locate the "Micro results" label cell in the (only) table, then grab
# the adjacent value cell in the same row - that is the cell whose long
# run of date-stamped micro-biology results needs to be summarised away.
$tbl = $d.Tables.Item(1)
$valueCell = $null
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $row = $tbl.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "*Micro results*") {
        $valueCell = $row.Cells.Item(2)
        break
    }
}

$cellRange = $valueCell.Range
# Exclude the trailing cell-mark character so we only target the
# paragraph content, not the end-of-cell marker.
$contentRange = $d.Range($cellRange.Start, $cellRange.End - 1)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newXml = "<w:p $w><w:r/></w:p>" +
          "<w:p $w><w:r><w:rPr>" +
          "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/>" +
          "<w:color w:val='0000FF'/>" +
          "<w:sz w:val='20'/>" +
          "</w:rPr></w:r></w:p>"

$contentRange.InsertXML($newXml)
